# fix: required fields on modal forms
#
# "Introducción a la Programación" (Sección 2) moves from the 11:00-13:00
# slot into the same 9:00-11:00 slot as Sección 1 (different room), freeing
# up the 11:00-12:00 slot on Mondays. "Algoritmos y Complejidad" (Sección 3)
# moves from Fridays 14:00-17:00 to Mondays 14:00-16:00 in a new room, and a
# new section ("asdf") is added on Fridays 15:00-17:00.

$wb  = $excel.ActiveWorkbook
$horario = $wb.Worksheets.Item("Horario")
$tabla   = $wb.Worksheets.Item("Tabla")

# ---- Horario (calendar grid) sheet ----

$introMerged = "Introducción a la Programación (1)`nReloj 102`nIntroducción a la Programación (2)`nReloj 103"
$algo        = "Algoritmos y Complejidad (3)`nReloj 103"
$introAsdf   = "Introducción a la Programación (asdf)`nCiencias 507"

# 9:00 and 10:00 Monday slots now show both Introducción sections
$horario.Cells.Item(2, 2).Value = $introMerged
$horario.Cells.Item(3, 2).Value = $introMerged

# 11:00 and 12:00 Monday slots are now free
$horario.Cells.Item(4, 2).Value = ""
$horario.Cells.Item(5, 2).Value = ""

# Algoritmos y Complejidad moves to Monday 14:00-16:00
$horario.Cells.Item(7, 2).Value = $algo
$horario.Cells.Item(8, 2).Value = $algo
$horario.Cells.Item(9, 2).Value = $algo

# Friday 14:00 slot is now free, 15:00-16:00 gets the new "asdf" section
$horario.Cells.Item(7, 6).Value = ""
$horario.Cells.Item(8, 6).Value = $introAsdf
$horario.Cells.Item(9, 6).Value = $introAsdf

# ---- Tabla (source data) sheet ----

# Sección 2: start moves 11:00 -> 9:00, room moves Reloj 102 -> Reloj 103
$tabla.Cells.Item(3, 3).Value = "9:00"
$tabla.Cells.Item(3, 4).Value = "11:00"
$tabla.Cells.Item(3, 5).Value = "Reloj 103"

# Sección 3: day moves Viernes -> Lunes, room moves Ciencias 507 -> Reloj 103
$tabla.Cells.Item(4, 2).Value = "Lunes"
$tabla.Cells.Item(4, 5).Value = "Reloj 103"

# New row: Sección asdf, Viernes 15:00-17:00, Ciencias 507
$tabla.Cells.Item(5, 1).Value = "Introducción a la Programación (Sección asdf)"
$tabla.Cells.Item(5, 2).Value = "Viernes"
$tabla.Cells.Item(5, 3).Value = "15:00"
$tabla.Cells.Item(5, 4).Value = "17:00"
$tabla.Cells.Item(5, 5).Value = "Ciencias 507"
